$wb = $excel.ActiveWorkbook

# ALC!row59
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 593.5
$ws.Range("I59").Value = 100
$ws.Range("J59").Value = 1251.5
$ws.Range("K59").Value = 300
$ws.Range("L59").Value = 3754.5
$ws.Range("M59").Value = 257
$ws.Range("N59").Value = -4868.5

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3027.762
$ws.Range("I113").Value = 2985.5715
$ws.Range("K113").Value = 2985.5715
$ws.Range("M113").Value = 268.4285

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11540.909
$ws.Range("I116").Value = 34601.668
$ws.Range("J116").Value = 2893.125
$ws.Range("K116").Value = 34601.668
$ws.Range("L116").Value = 2893.125
$ws.Range("M116").Value = -31159.668
$ws.Range("N116").Value = -9777.125

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 29413894
$ws.Range("I132").Value = 40000908
$ws.Range("J132").Value = 5522.4443
$ws.Range("K132").Value = 120002724
$ws.Range("L132").Value = 16567.3329
$ws.Range("M132").Value = -120000194
$ws.Range("N132").Value = -21627.3329

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2508.0657
$ws.Range("I138").Value = 1398.3585
$ws.Range("J138").Value = 5065.2173
$ws.Range("K138").Value = 4195.0755
$ws.Range("L138").Value = 15195.6519
$ws.Range("M138").Value = 944.9245000000001
$ws.Range("N138").Value = -25475.6519

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 331068.06
$ws.Range("I32").Value = 2185.8481
$ws.Range("J32").Value = 3578780
$ws.Range("K32").Value = 2185.8481
$ws.Range("L32").Value = 3578780
$ws.Range("M32").Value = -1898.8481
$ws.Range("N32").Value = -3579354

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1216.4318
$ws.Range("I61").Value = 822.7406999999999
$ws.Range("J61").Value = 1841.7059
$ws.Range("K61").Value = 822.7406999999999
$ws.Range("L61").Value = 1841.7059
$ws.Range("M61").Value = -610.7406999999999
$ws.Range("N61").Value = -2265.7059

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 6658.7144
$ws.Range("I102").Value = 3300
$ws.Range("K102").Value = 3300
$ws.Range("M102").Value = -1678

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 26333942
$ws.Range("I122").Value = 28590952
$ws.Range("J122").Value = 2180.6667
$ws.Range("K122").Value = 85772856
$ws.Range("L122").Value = 6542.000100000001
$ws.Range("M122").Value = -85770406
$ws.Range("N122").Value = -11442.0001

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1860.7097
$ws.Range("I132").Value = 975.0476
$ws.Range("J132").Value = 3720.6
$ws.Range("K132").Value = 2925.1428
$ws.Range("L132").Value = 11161.8
$ws.Range("M132").Value = -395.1428000000001
$ws.Range("N132").Value = -16221.8

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1216.4318
$ws.Range("I136").Value = 822.7406999999999
$ws.Range("J136").Value = 1841.7059
$ws.Range("K136").Value = 2468.2221
$ws.Range("L136").Value = 5525.1177
$ws.Range("M136").Value = 81.77790000000005
$ws.Range("N136").Value = -10625.1177

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 742.79486
$ws.Range("I94").Value = 723.40625
$ws.Range("J94").Value = 831.4286
$ws.Range("K94").Value = 723.40625
$ws.Range("L94").Value = 831.4286
$ws.Range("M94").Value = -272.40625
$ws.Range("N94").Value = -1733.4286

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1518.1111
$ws.Range("I105").Value = 1451.579
$ws.Range("J105").Value = 1592.4706
$ws.Range("K105").Value = 1451.579
$ws.Range("L105").Value = 1592.4706
$ws.Range("M105").Value = 295.421
$ws.Range("N105").Value = -5086.470600000001

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 891.1111
$ws.Range("I107").Value = 744.5
$ws.Range("J107").Value = 1310
$ws.Range("K107").Value = 744.5
$ws.Range("L107").Value = 1310
$ws.Range("M107").Value = 1175.5
$ws.Range("N107").Value = -5150

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7965.619
$ws.Range("I134").Value = 943.2222
$ws.Range("J134").Value = 50100
$ws.Range("K134").Value = 2829.6666
$ws.Range("L134").Value = 150300
$ws.Range("M134").Value = -294.6666
$ws.Range("N134").Value = -155370

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12197469
$ws.Range("I31").Value = 20001976
$ws.Range("J31").Value = 2927.375
$ws.Range("K31").Value = 20001976
$ws.Range("L31").Value = 2927.375
$ws.Range("M31").Value = -20001681
$ws.Range("N31").Value = -3517.375

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 12197469
$ws.Range("I34").Value = 20001976
$ws.Range("J34").Value = 2927.375
$ws.Range("K34").Value = 20001976
$ws.Range("L34").Value = 2927.375
$ws.Range("M34").Value = -20001774
$ws.Range("N34").Value = -3331.375

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 793.4894
$ws.Range("I58").Value = 607.56757
$ws.Range("J58").Value = 1481.4
$ws.Range("K58").Value = 607.56757
$ws.Range("L58").Value = 1481.4
$ws.Range("M58").Value = -404.56757
$ws.Range("N58").Value = -1887.4

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1566.6111
$ws.Range("I99").Value = 1538.4615
$ws.Range("J99").Value = 1639.8
$ws.Range("K99").Value = 1538.4615
$ws.Range("L99").Value = 1639.8
$ws.Range("M99").Value = -40.46149999999989
$ws.Range("N99").Value = -4635.8

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1451.5588
$ws.Range("I122").Value = 1459.9584
$ws.Range("J122").Value = 1431.4
$ws.Range("K122").Value = 4379.8752
$ws.Range("L122").Value = 4294.200000000001
$ws.Range("M122").Value = -1929.8752
$ws.Range("N122").Value = -9194.200000000001

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1566.6111
$ws.Range("I126").Value = 1538.4615
$ws.Range("J126").Value = 1639.8
$ws.Range("K126").Value = 4615.3845
$ws.Range("L126").Value = 4919.4
$ws.Range("M126").Value = -2145.3845
$ws.Range("N126").Value = -9859.4

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 45359.39
$ws.Range("I132").Value = 904.7857
$ws.Range("K132").Value = 2714.3571
$ws.Range("M132").Value = -184.3571000000002

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 23220.133
$ws.Range("I134").Value = 27365.947
$ws.Range("J134").Value = 714.2857
$ws.Range("K134").Value = 82097.841
$ws.Range("L134").Value = 2142.8571
$ws.Range("M134").Value = -79562.841
$ws.Range("N134").Value = -7212.8571

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 793.4894
$ws.Range("I136").Value = 607.56757
$ws.Range("J136").Value = 1481.4
$ws.Range("K136").Value = 1822.70271
$ws.Range("L136").Value = 4444.200000000001
$ws.Range("M136").Value = 727.29729
$ws.Range("N136").Value = -9544.200000000001

# CUL!row117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 950
$ws.Range("I117").Value = 950
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 2850
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 592
$ws.Range("N117").ClearContents()

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 64822.344
$ws.Range("I132").Value = 57511.277
$ws.Range("J132").Value = 74222.28999999999
$ws.Range("K132").Value = 172533.831
$ws.Range("L132").Value = 222666.87
$ws.Range("M132").Value = -170003.831
$ws.Range("N132").Value = -227726.87

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5757.154
$ws.Range("I136").Value = 3283.111
$ws.Range("J136").Value = 11323.75
$ws.Range("K136").Value = 9849.332999999999
$ws.Range("L136").Value = 33971.25
$ws.Range("M136").Value = -7299.332999999999
$ws.Range("N136").Value = -39071.25

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2006.5172
$ws.Range("I122").Value = 1483.0555
$ws.Range("J122").Value = 2863.0908
$ws.Range("K122").Value = 4449.166499999999
$ws.Range("L122").Value = 8589.2724
$ws.Range("M122").Value = -1999.166499999999
$ws.Range("N122").Value = -13489.2724

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4345.9653
$ws.Range("I136").Value = 5246.9546
$ws.Range("J136").Value = 1514.2858
$ws.Range("K136").Value = 15740.8638
$ws.Range("L136").Value = 4542.857400000001
$ws.Range("M136").Value = -13190.8638
$ws.Range("N136").Value = -9642.857400000001
